# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
# - Inserts a new detail row (period 2509) below the existing 2507/2508 rows
# - Updates the "Cant. Periodos" count and the "VALOR MORA" total accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last existing detail row (17) and insert it as the new row 18,
# pushing the blank rows and the signature block down by one row - this is
# what Excel's "Insert Copied Cells" does and keeps the table's formatting.
$ws.Rows.Item(17).Copy()
$ws.Rows.Item(18).Insert()

# The trailing (empty) cells don't carry the border formatting over on
# insert, so re-apply just the formatting from the row above.
$ws.Range("H17:J17").Copy()
$ws.Range("H18").PasteSpecial(-4122)

# Fill the new row with the new period's data (same worker, new period 2509)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1002474995"
$ws.Range("D18").Value = "YELIANA LUCIA VASQUEZ OROZCO"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Keep the "Periodo Mora" column centered, same as the rest of the table
$ws.Range("E16:E18").HorizontalAlignment = -4108

# Update the totals: one more period (3) and the accumulated Valor Mora
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = 170820
